# IssueTracking.xlsx - "Began to add MessageConfiguration class"
#
# - Bold the header row (A1:C1)
# - Move the status-legend cells from column I (row-per-status) into a
#   single header row, columns E:H, colored with the Good/Bad/Neutral
#   named cell styles
# - Drop the old per-row "status dot" cells in column H (rows 2-4) and the
#   old legend text in column I (rows 1-4)
# - Issue #3 ("Create MessageConfiguration class") is now closed -> color
#   its row with the "Good" style
# - Add four new backlog rows (11-14) for the new MessageConfiguration work
# - Freeze the header row and leave the selection on B15
# - Force portrait page orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: bold A1:C1 ---
$ws.Range("A1:C1").Font.Bold = $true

# --- Status legend: rebuild on row 1, columns E:H ---
$ws.Range("I1").Clear() | Out-Null

$ws.Range("E1").Value = "open"

$ws.Range("F1").Value = "closed"
$ws.Range("F1").Style = "Good"

$ws.Range("G1").Value = "urgent"
$ws.Range("G1").Style = "Bad"

$ws.Range("H1").Value = "deleted"
$ws.Range("H1").Style = "Neutral"

# --- Remove the old per-row status cells (H2:H4 swatches, I2:I4 labels) ---
$ws.Range("H2:H4").Clear() | Out-Null
$ws.Range("I2:I4").Clear() | Out-Null

# --- Issue #3 is now closed -> highlight its row ---
$ws.Range("A4:C4").Style = "Good"

# --- New backlog rows ---
$ws.Range("B11").Value = "Citation for MCOS paper here"
$ws.Range("C11").Value = "main"

$ws.Range("B12").Value = "WeightingScheme to be defined"
$ws.Range("C12").Value = "config"

$ws.Range("B13").Value = "Create a method that can pull out information for a given OMT"
$ws.Range("C13").Value = "MessageConfiguration"

$ws.Range("B14").Value = "Build a constructor for MessageConfiguration"
$ws.Range("C14").Value = "MessageConfiguration"

# --- Freeze header row, leave selection on B15 ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B15").Select() | Out-Null

# --- Page setup ---
$ws.PageSetup.Orientation = 1
